$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the time value for E6 (it will propagate through the shared
# formulas in E13 and E14 via recalculation)
$ws.Range("E6").Value = 360

# Force a full recalculation so dependent formula cells (E13, E14) update
$excel.Calculate()

# Update the selected/active cell shown when the workbook is reopened
$ws.Range("E7").Select()
